# Update ligand/receptor expression & derived-specificity values for the
# Thbs1-Ptprj LR-pairs sheet (rows 2-26) to reflect the new TPM-based
# recomputation. Two contiguous column blocks change:
#   G:J  - Ligand average/total expression value + derived specificity
#   M:T  - Receptor average/total expression value + derived specificity,
#          and edge average/total expression weight + derived specificity
# Columns A-F, K, L are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for G2:J26 (one row per array element, in sheet-row order)
$dataGJ = @(
    @(16.71895933333333,50.156878,0.02912144738161902,0.03059269312988411),
    @(16.71895933333333,50.156878,0.02912144738161902,0.03059269312988411),
    @(16.71895933333333,50.156878,0.02912144738161902,0.03059269312988411),
    @(16.71895933333333,50.156878,0.02912144738161902,0.03059269312988411),
    @(16.71895933333333,50.156878,0.02912144738161902,0.03059269312988411),
    @(155.500389,466.501167,0.2708539632042961,0.2845377865576845),
    @(155.500389,466.501167,0.2708539632042961,0.2845377865576845),
    @(155.500389,466.501167,0.2708539632042961,0.2845377865576845),
    @(155.500389,466.501167,0.2708539632042961,0.2845377865576845),
    @(155.500389,466.501167,0.2708539632042961,0.2845377865576845),
    @(194.8548433333333,584.56453,0.3394024086099587,0.3565493705749576),
    @(194.8548433333333,584.56453,0.3394024086099587,0.3565493705749576),
    @(194.8548433333333,584.56453,0.3394024086099587,0.3565493705749576),
    @(194.8548433333333,584.56453,0.3394024086099587,0.3565493705749576),
    @(194.8548433333333,584.56453,0.3394024086099587,0.3565493705749576),
    @(82.82950199999999,165.659004,0.1442742299952585,0.1010420758958371),
    @(82.82950199999999,165.659004,0.1442742299952585,0.1010420758958371),
    @(82.82950199999999,165.659004,0.1442742299952585,0.1010420758958371),
    @(82.82950199999999,165.659004,0.1442742299952585,0.1010420758958371),
    @(82.82950199999999,165.659004,0.1442742299952585,0.1010420758958371),
    @(124.2078576666667,372.623573,0.2163479508088675,0.2272780738416368),
    @(124.2078576666667,372.623573,0.2163479508088675,0.2272780738416368),
    @(124.2078576666667,372.623573,0.2163479508088675,0.2272780738416368),
    @(124.2078576666667,372.623573,0.2163479508088675,0.2272780738416368),
    @(124.2078576666667,372.623573,0.2163479508088675,0.2272780738416368)
)

# New values for M2:T26 (one row per array element, in sheet-row order)
$dataMT = @(
    @(1.135186,3.405558,0.006415563145489937,0.00646579730555003,18.97912856976933,170.812157127924,0.0001868304845648394,0.0001978061528087236),
    @(4.728754666666667,14.186264,0.02672480471352731,0.02693406118674866,79.05985696931023,711.5387127237921,0.0007782649942490295,0.0008239854686277241),
    @(84.55360633333333,253.660819,0.4778591355164685,0.4816007949398642,1413.648305773676,12722.83475196308,0.01391594967076879,0.01473346533070351),
    @(4.124113,8.248226,0.02330764066032874,0.01566009371925767,68.95087753307135,413.705265198428,0.0006787522310794473,0.0004790844415384754),
    @(82.400874,247.202622,0.4656928559641855,0.4693392528485795,1377.656861437124,12398.91175293412,0.01356165000095691,0.01435835173620568),
    @(1.135186,3.405558,0.006415563145489937,0.00646579730555003,176.521864587354,1588.696781286186,0.001737680704143369,0.001839763653651846),
    @(4.728754666666667,14.186264,0.02672480471352731,0.02693406118674866,735.3231901522321,6617.908711370089,0.007238519272519725,0.007663758153086706),
    @(84.55360633333333,253.660819,0.4778591355164685,0.4816007949398642,13148.1186761862,118333.0680856758,0.1294300407080143,0.1370336241966103),
    @(4.124113,8.248226,0.02330764066032874,0.01566009371925767,641.3011757799571,3847.807054679743,0.006312966845791635,0.004455888404163475),
    @(82.400874,247.202622,0.4656928559641855,0.4693392528485795,12813.36796093999,115320.3116484599,0.126134755673827,0.1335447521501723),
    @(1.135186,3.405558,0.006415563145489937,0.00646579730555003,221.1964901841933,1990.76841165774,0.002177457584168568,0.00230537595955912),
    @(4.728754666666667,14.186264,0.02672480471352731,0.02693406118674866,921.4207497351023,8292.78674761592,0.009070463089401948,0.00960332256316263),
    @(84.55360633333333,253.660819,0.4778591355164685,0.4816007949398642,16475.67971535001,148281.1174381501,0.1621865415705621,0.1717144603042078),
    @(4.124113,8.248226,0.02330764066032874,0.01566009371925767,803.6033925039634,4821.620355023781,0.007910669379130984,0.005583596558746169),
    @(82.400874,247.202622,0.4656928559641855,0.4693392528485795,16056.20939379974,144505.8845441977,0.1580572769866952,0.1673426151892819),
    @(1.135186,3.405558,0.006415563145489937,0.00646579730555003,94.026891057372,564.1613463442319,0.0009256004328015192,0.0006533175820744854),
    @(4.728754666666667,14.186264,0.02672480471352731,0.02693406118674866,391.680394120176,2350.082364721056,0.003855700621817808,0.00272147345461458),
    @(84.55360633333333,253.660819,0.4778591355164685,0.4816007949398642,7003.533104894045,42021.19862936427,0.06894275882283837,0.04866194407380926),
    @(4.124113,8.248226,0.02330764066032874,0.01566009371925767,341.598225981726,1366.392903926904,0.003362691909275107,0.001582328378117156),
    @(82.400874,247.202622,0.4656928559641855,0.4693392528485795,6825.223357784747,40951.34014670848,0.06718747820852569,0.04742301240722167),
    @(1.135186,3.405558,0.006415563145489937,0.00646579730555003,140.9990211131927,1268.991190018734,0.00138799393981164,0.001469533957455856),
    @(4.728754666666667,14.186264,0.02672480471352731,0.02693406118674866,587.3484865779191,5286.136379201272,0.005781856735538799,0.006121521547257027),
    @(84.55360633333333,253.660819,0.4778591355164685,0.4816007949398642,10502.22230065403,94520.00070588628,0.1033838447442849,0.1094573010345335),
    @(4.124113,8.248226,0.02330764066032874,0.01566009371925767,512.2472405052497,3073.483443031498,0.005042560295051564,0.003559195936692398),
    @(82.400874,247.202622,0.4656928559641855,0.4693392528485795,10234.83602940093,92113.52426460841,0.1007516950941807,0.1066705213656981)
)

$nRows = $dataGJ.Count

$arrGJ = New-Object 'object[,]' $nRows,4
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $arrGJ[$i,$j] = $dataGJ[$i][$j]
    }
}

$arrMT = New-Object 'object[,]' $nRows,8
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt 8; $j++) {
        $arrMT[$i,$j] = $dataMT[$i][$j]
    }
}

$ws.Range("G2:J26").Value2 = $arrGJ
$ws.Range("M2:T26").Value2 = $arrMT
